$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("大位科技", "利欧股份", "利欧股份"),
    @("英维克",   "协鑫集成", "巨力索具"),
    @("首都在线", "巨力索具", "双良节能"),
    @("利欧股份", "蓝色光标", "协鑫集成"),
    @("掌阅科技", "大位科技", "嘉美包装"),
    @("巨力索具", "特变电工", "大位科技"),
    @("协鑫集成", "英维克",   "航天发展"),
    @("蓝色光标", "掌阅科技", "蓝色光标"),
    @("网宿科技", "双良节能", "首都在线"),
    @("昆仑万维", "格林美",   "特发信息"),
    @("双良节能", "首都在线", "天际股份"),
    @("特发信息", "特发信息", "雅博股份"),
    @("航天发展", "网宿科技", "英维克"),
    @("特变电工", "昆仑万维", "掌阅科技"),
    @("天际股份", "天孚通信", "网宿科技"),
    @("优刻得-W", "汉缆股份", "博纳影业"),
    @("格林美",   "光线传媒", "杭电股份"),
    @("浙文互联", "航天发展", "光线传媒"),
    @("光线传媒", "天际股份", "锋龙股份"),
    @("天孚通信", "浙文互联", "浙文互联")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
